# Rename the header cells: "FactorOne"/"FactorTwo" -> "Number1"/"Number2"
# (the "Result" header stays the same), and move the active selection to B2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Number1"
$ws.Range("B1").Value = "Number2"
$ws.Range("C1").Value = "Result"

$ws.Range("B2").Select()
